# Atom_TMA.xlsx : E-core TMA update 3.51 / Atom TMA update 1.2
# - Bump version number on "Formulas" sheet (C1: 1.1 -> 1.2)
# - Remove trailing blank row (row 79) on "Formulas" sheet
# - Spelling / trailing-space clean-up of several metric descriptions
# - Add a new "1.2" row to the "Change Log" sheet documenting this update

$wb = $excel.ActiveWorkbook
$formulas = $wb.Worksheets.Item("Formulas")
$changeLog = $wb.Worksheets.Item("Change Log")

# ---- Formulas sheet -------------------------------------------------

# Version bump
$formulas.Range("C1").Value = 1.2

# Spelling fixes / trailing-space removal on metric descriptions
$formulas.Range("J20").Value = "Counts the number of cycles due to backend bound stalls that are core execution bound and not attributed to outstanding demand load or store stalls."
$formulas.Range("J21").Value = "Counts the number of cycles the core is stalled due to stores or loads."
$formulas.Range("J26").Value = "Counts the total number of issue slots  that were not consumed by the backend due to backend stalls.  Note that UOPS must be available for consumption in order for this event to count.  If a uop is not available (IQ is empty), this event will not count.  All of these subevents count backend stalls, in slots, due to a resource limitation.   These are not cycle based events and therefore can not be precisely added or subtracted from the Backend_Bound subevents which are cycle based.  These subevents are supplementary to Backend_Bound and can be used to analyze results from a resource perspective at allocation."
$formulas.Range("J27").Value = "Counts the total number of issue slots  that were not consumed by the backend due to backend stalls.  Note that uops must be available for consumption in order for this event to count.  If a uop is not available (IQ is empty), this event will not count."
$formulas.Range("J35").Value = "Counts the number of issue slots  that result in retirement slots."
$formulas.Range("J36").Value = "Counts the number of uops that are not from the microsequencer."
$formulas.Range("J50").Value = "Instructions per Branch (lower number means higher occurence rate)"
$formulas.Range("J51").Value = "Instruction per (near) call (lower number means higher occurence rate)"

# Remove the trailing empty row at the bottom of the sheet
$formulas.Rows.Item(79).Delete()

# ---- Change Log sheet ------------------------------------------------

# Make room for the new entry at the top of the log (below the header row)
$changeLog.Rows.Item(2).Insert()

$changeLog.Range("A2").Value = 1.2
$changeLog.Range("B2").Value = "Spelling fixes for occurrence, number; Removed trailing spaces on several descriptions"
$changeLog.Range("C2").Value = 45170
